$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D4: new "Iniciar sesión" label (adds a new shared string)
$ws.Range("D4").Value = "Iniciar sesión"

# Row 5: fill in B5 (hyperlinked email, same style as B2:B4) and C5 (numeric value)
$ws.Range("B5").Value = "jisola.tsoft@gmail.com"
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:jisola.tsoft@gmail.com")
$ws.Range("B5").Style = "Hyperlink"
$ws.Range("C5").Value = 12061990

# Move the active selection from D4 to C4
$ws.Range("C4").Select()
